# Auto-generated edit script: applies the Alpha_Profits value updates
# described in the commit diff, per-sheet / per-row, via Excel COM.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 744332.75
$ws.Range("I8").Value = 837249.5
$ws.Range("J8").Value = 999
$ws.Range("K8").Value = 2511748.5
$ws.Range("L8").Value = 2997
$ws.Range("M8").Value = -2511609.5
$ws.Range("N8").Value = -3275

$ws.Range("H43").Value = 1677992
$ws.Range("I43").Value = 2235656
$ws.Range("K43").Value = 2235656
$ws.Range("M43").Value = -2235587

$ws.Range("H51").Value = 16219.873
$ws.Range("I51").Value = 4804.523
$ws.Range("J51").Value = 61881.273
$ws.Range("K51").Value = 4804.523
$ws.Range("L51").Value = 61881.273
$ws.Range("M51").Value = -4320.523
$ws.Range("N51").Value = -62849.273

$ws.Range("H69").Value = 4960.7856
$ws.Range("I69").Value = 4960.7856
$ws.Range("K69").Value = 14882.3568
$ws.Range("M69").Value = -14008.3568

$ws.Range("H72").Value = 4960.7856
$ws.Range("I72").Value = 4960.7856
$ws.Range("K72").Value = 44647.0704
$ws.Range("M72").Value = -40279.0704

$ws.Range("H87").Value = 30000
$ws.Range("J87").Value = 30000
$ws.Range("L87").Value = 30000
$ws.Range("N87").Value = -32496

$ws.Range("H88").Value = 1503.8
$ws.Range("J88").Value = 957.4
$ws.Range("L88").Value = 957.4
$ws.Range("N88").Value = -1769.4

$ws.Range("H90").Value = 30000
$ws.Range("J90").Value = 30000
$ws.Range("L90").Value = 90000
$ws.Range("N90").Value = -102480

$ws.Range("H91").Value = 1503.8
$ws.Range("J91").Value = 957.4
$ws.Range("L91").Value = 957.4
$ws.Range("N91").Value = -3765.4

$ws.Range("H93").Value = 87650.25
$ws.Range("J93").Value = 87650.25
$ws.Range("L93").Value = 87650.25
$ws.Range("N93").Value = -92642.25

$ws.Range("H94").Value = 1527.6
$ws.Range("I94").Value = 1330.6666
$ws.Range("K94").Value = 1330.6666
$ws.Range("M94").Value = -879.6666

$ws.Range("H138").Value = 3385.6086
$ws.Range("I138").Value = 2191.6
$ws.Range("J138").Value = 3454.23
$ws.Range("K138").Value = 6574.799999999999
$ws.Range("L138").Value = 10362.69
$ws.Range("M138").Value = -1434.799999999999
$ws.Range("N138").Value = -20642.69

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 42999
$ws.Range("J44").Value = 42999
$ws.Range("L44").Value = 42999
$ws.Range("N44").Value = -43975

$ws.Range("H55").Value = 39000
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 39000
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 39000
$ws.Range("N55").Value = -39630
$ws.Range("M55").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 12795.6455
$ws.Range("I86").Value = 10357.529
$ws.Range("K86").Value = 10357.529
$ws.Range("M86").Value = -9234.529

$ws.Range("H89").Value = 12795.6455
$ws.Range("I89").Value = 10357.529
$ws.Range("K89").Value = 51787.645
$ws.Range("M89").Value = -46171.645

$ws.Range("H105").Value = 106541.055
$ws.Range("I105").Value = 1392.8125
$ws.Range("K105").Value = 1392.8125
$ws.Range("M105").Value = 354.1875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1098.7368
$ws.Range("I16").Value = 1061.6666
$ws.Range("J16").Value = 1237.75
$ws.Range("K16").Value = 1061.6666
$ws.Range("L16").Value = 1237.75
$ws.Range("M16").Value = -774.6666
$ws.Range("N16").Value = -1811.75

$ws.Range("H43").Value = 38042.855
$ws.Range("J43").Value = 38042.855
$ws.Range("L43").Value = 38042.855
$ws.Range("N43").Value = -38410.855

$ws.Range("H54").Value = 22690.154
$ws.Range("J54").Value = 23331
$ws.Range("L54").Value = 23331
$ws.Range("N54").Value = -24647

$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()

$ws.Range("H95").Value = 15401.417
$ws.Range("J95").Value = 15401.417
$ws.Range("L95").Value = 15401.417
$ws.Range("N95").Value = -20893.417

$ws.Range("H101").Value = 38042.855
$ws.Range("J101").Value = 38042.855
$ws.Range("L101").Value = 38042.855
$ws.Range("N101").Value = -44532.855

$ws.Range("H105").Value = 2055
$ws.Range("I105").Value = 2067.7693
$ws.Range("J105").Value = 1999.6666
$ws.Range("K105").Value = 2067.7693
$ws.Range("L105").Value = 1999.6666
$ws.Range("M105").Value = -320.7692999999999
$ws.Range("N105").Value = -5493.6666

$ws.Range("H113").Value = 1098.7368
$ws.Range("I113").Value = 1061.6666
$ws.Range("J113").Value = 1237.75
$ws.Range("K113").Value = 1061.6666
$ws.Range("L113").Value = 1237.75
$ws.Range("M113").Value = 1108.3334
$ws.Range("N113").Value = -5577.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 5422
$ws.Range("I56").Value = 5422
$ws.Range("K56").Value = 5422
$ws.Range("M56").Value = -4892

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4763.2856
$ws.Range("I113").Value = 4757.6
$ws.Range("K113").Value = 4757.6
$ws.Range("M113").Value = -2587.6

$ws.Range("H126").Value = 16413.143
$ws.Range("I126").Value = 15960.4
$ws.Range("J126").Value = 16664.666
$ws.Range("K126").Value = 47881.2
$ws.Range("L126").Value = 49993.99800000001
$ws.Range("M126").Value = -45411.2
$ws.Range("N126").Value = -54933.99800000001

$ws.Range("H130").Value = 78000
$ws.Range("J130").Value = 78000
$ws.Range("L130").Value = 78000
$ws.Range("N130").Value = -88040

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3281.0908
$ws.Range("I82").Value = 3911.625
$ws.Range("J82").Value = 1599.6666
$ws.Range("K82").Value = 3911.625
$ws.Range("L82").Value = 1599.6666
$ws.Range("M82").Value = -3550.625
$ws.Range("N82").Value = -2321.6666

$ws.Range("H85").Value = 3281.0908
$ws.Range("I85").Value = 3911.625
$ws.Range("J85").Value = 1599.6666
$ws.Range("K85").Value = 3911.625
$ws.Range("L85").Value = 1599.6666
$ws.Range("M85").Value = -2663.625
$ws.Range("N85").Value = -4095.6666

$ws.Range("H122").Value = 6800.8
$ws.Range("I122").Value = 5565.1816
$ws.Range("J122").Value = 10198.75
$ws.Range("K122").Value = 16695.5448
$ws.Range("L122").Value = 30596.25
$ws.Range("M122").Value = -14245.5448
$ws.Range("N122").Value = -35496.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3480
$ws.Range("I62").Value = 3450
$ws.Range("K62").Value = 3450
$ws.Range("M62").Value = -2826

$ws.Range("H65").Value = 3480
$ws.Range("I65").Value = 3450
$ws.Range("K65").Value = 17250
$ws.Range("M65").Value = -14130

$ws.Range("H80").Value = 21558
$ws.Range("J80").Value = 23150
$ws.Range("L80").Value = 23150
$ws.Range("N80").Value = -25146

$ws.Range("H81").Value = 8112.125
$ws.Range("J81").Value = 10466
$ws.Range("L81").Value = 20932
$ws.Range("N81").Value = -23054

$ws.Range("H82").Value = 29900.5
$ws.Range("J82").Value = 29900.5
$ws.Range("L82").Value = 29900.5
$ws.Range("N82").Value = -30666.5

$ws.Range("H83").Value = 21558
$ws.Range("J83").Value = 23150
$ws.Range("L83").Value = 69450
$ws.Range("N83").Value = -79434

$ws.Range("H84").Value = 8112.125
$ws.Range("J84").Value = 10466
$ws.Range("L84").Value = 104660
$ws.Range("N84").Value = -115268

$ws.Range("H85").Value = 29900.5
$ws.Range("J85").Value = 29900.5
$ws.Range("L85").Value = 29900.5
$ws.Range("N85").Value = -32552.5

$ws.Range("H96").Value = 2610.6667
$ws.Range("I96").Value = 2563.7144
$ws.Range("K96").Value = 2563.7144
$ws.Range("M96").Value = -1190.7144
